$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '67.951.35'
$ws.Cells.Item(2, 5).Value = '  +0.14%  '
$ws.Cells.Item(3, 4).Value = '3.790.58'
$ws.Cells.Item(3, 5).Value = '  -1.03%  '
$ws.Cells.Item(4, 5).Value = '  +0.09%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '602.00'
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = '  -0.44%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '163.21'
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = '  -2.35%  '
$ws.Cells.Item(7, 4).Value = '3.787.90'
$ws.Cells.Item(7, 5).Value = '  -1.05%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.998'
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).Value = '  -0.25%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.515'
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).Value = '  -1.09%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.158'
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value = '  -2.23%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '6.89'
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value = '  +9.32%  '
$ws.Cells.Item(12, 5).Value = '  -1.43%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.0000247'
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).Value = '  -2.85%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '35.07'
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).Value = '  -2.72%  '
$ws.Cells.Item(15, 4).Value = '4.423.87'
$ws.Cells.Item(15, 5).Value = '  -0.95%  '
$ws.Cells.Item(16, 4).Value = '3.757.00'
$ws.Cells.Item(16, 5).Value = '  -2.44%  '
$ws.Cells.Item(17, 4).Value = '67.907.25'
$ws.Cells.Item(17, 5).Value = '  +0.06%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '18.20'
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).Value = '  -1.63%  '
$ws.Cells.Item(19, 5).Value = '  +2.00%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '7.00'
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = '  -1.35%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '458.94'
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).Value = '  -1.18%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '9.45'
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = '  -4.95%  '
$ws.Cells.Item(23, 5).Value = '  -1.49%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '83.19'
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value = '  -0.38%  '
$ws.Cells.Item(25, 5).Value = '  -4.67%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '11.88'
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).Value = '  -1.67%  '
$ws.Cells.Item(27, 5).Value = '  -1.64%  '
$ws.Cells.Item(28, 5).Value = '  -0.26%  '
$ws.Cells.Item(29, 5).Value = '  -1.67%  '
$ws.Cells.Item(30, 4).Value = '3.937.21'
$ws.Cells.Item(30, 5).Value = '  -0.92%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '7.22'
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '2.60'
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(32, 5).Value = '  -7.12%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '2.19'
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = '  -2.09%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '29.01'
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = '  -2.60%  '
$ws.Cells.Item(35, 5).Value = '  +0.00%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '8.92'
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).Value = '  -2.04%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.0992'
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value = '  -1.21%  '
$ws.Cells.Item(38, 5).Value = '  +5.31%  '
$ws.Cells.Item(39, 5).Value = '  -0.14%  '
$ws.Cells.Item(40, 2).Value = 'dogwifhat'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '3.20'
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).Value = '  -4.84%  '
$ws.Cells.Item(41, 2).Value = 'Mantle'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.979'
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).Value = '  -2.01%  '
$ws.Cells.Item(42, 5).Value = '  +0.07%  '
$ws.Cells.Item(43, 5).Value = '  +0.10%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '43.76'
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value = '  +1.01%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '47.14'
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = '  -2.01%  '
$ws.Cells.Item(46, 5).Value = '  +2.59%  '
$ws.Cells.Item(47, 5).Value = '  -2.37%  '
$ws.Cells.Item(48, 5).Value = '  -2.85%  '
$ws.Cells.Item(49, 5).Value = '  -1.07%  '
$ws.Cells.Item(50, 5).Value = '  -0.96%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '26.67'
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value = '  -5.88%  '
